Write-Output "hello"
Write-Output 42
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$v = $ws.Range("N3")
Write-Output $v.GetType()
